# Weekly update: insert a new week's worth of data (4 quality rows, dated
# 2023-03-10 / serial 44995) at the top of the "Pepino dulce" price table,
# pushing the existing history down by four rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows above the current row 279 (the first data row of the
# block), shifting all the existing rows (old 279-316) down to (283-320).
$ws.Range("A279:A282").EntireRow.Insert()

# Common values shared by every row of this market/product block.
$mercadoId = 6
$mercado   = "Mercado Mayorista Lo Valledor de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$categoriaId = 100112043
$categoria = "Pepino dulce"
$variedad  = "Cultivar IV Región"
$unidad    = "`$/bandeja 18 kilos"
$origen    = "Provincia de Limarí"
$kgUnidades = 18
$clasificacion = "Hortaliza"
$fecha = 44995

# New rows: [row, calidad, volumen, precioMin, precioMax, precioProm, precioKg]
$newRows = @(
    @(279, "Especial", 260, 15000, 15000, 15000, 833),
    @(280, "Primera",  420, 14000, 14000, 14000, 778),
    @(281, "Segunda",  115, 12000, 12000, 12000, 667),
    @(282, "Tercera",   70, 10000, 10000, 10000, 556)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value  = $mercadoId
    $ws.Cells.Item($r, 2).Value  = $mercado
    $ws.Cells.Item($r, 3).Value  = $region
    $ws.Cells.Item($r, 4).Value  = $fecha
    $ws.Cells.Item($r, 5).Value  = $codreg
    $ws.Cells.Item($r, 6).Value  = $categoriaId
    $ws.Cells.Item($r, 7).Value  = $categoria
    $ws.Cells.Item($r, 8).Value  = $variedad
    $ws.Cells.Item($r, 9).Value  = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
    $ws.Cells.Item($r, 11).Value = $row[3]
    $ws.Cells.Item($r, 12).Value = $row[4]
    $ws.Cells.Item($r, 13).Value = $row[5]
    $ws.Cells.Item($r, 14).Value = $unidad
    $ws.Cells.Item($r, 15).Value = $origen
    $ws.Cells.Item($r, 16).Value = $row[6]
    $ws.Cells.Item($r, 17).Value = $kgUnidades
    $ws.Cells.Item($r, 18).Value = $clasificacion
}
